# Updated cryptos list on Mon Sep 23 18:49:23 UTC 2024 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns with latest scrape; a few
# coins (Bittensor/Binance-PegBSC-USD, VeChain/Stellar) swapped rank
# positions, so their Coin/Link/Price/Volume cells are rewritten too.
#
# NumberFormat is forced to "@" (Text) immediately before assigning any
# Price value that would otherwise be auto-parsed by Excel as a number
# (and so lose a significant trailing zero, e.g. "11.30" -> 11.3) - this
# keeps the cell's stored type/text identical to the original inline
# string, matching values such as "0.999", "11.30", "0.0000146".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.525.33"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.684.59"
$ws.Range("E3").Value = "  +4.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.05"
$ws.Range("E5").Value = "  +4.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.11"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "2.684.51"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.62"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.359"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.53"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "3.161.66"
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").Value = "63.391.30"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000146"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "2.687.70"
$ws.Range("E18").Value = "  +4.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.48"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.10"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.43"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.58"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E25").Value = "  +3.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.55"
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.64"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "543.75"
$ws.Range("E29").Value = "  +18.06%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  +7.20%  "
$ws.Range("E33").Value = "  +7.45%  "
$ws.Range("D34").Value = "0.0₃0812"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "172.55"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("E36").Value = "  +12.85%  "
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.35"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("E40").Value = "  +8.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.31"
$ws.Range("E41").Value = "  +11.71%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.26"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0568"
$ws.Range("E45").Value = "  +5.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.637"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0241"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0966"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.03"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("E50").Value = "  +4.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.30"
$ws.Range("E51").Value = "  -0.92%  "
